$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project 1")

# Fill in row 5 of the Project1 table with the new "Fibonaaci" DP entry.
$ws.Range("B5").Value = "DP"
$ws.Range("C5").Value = "Fibonaaci"
$ws.Range("E5").Value = "O(n)"
$ws.Range("F5").Value = "Store the computed sub-problem in hashTable. If the value is found in the hashTable, return else store it."

# Update the active selection to match the author's saved state.
$ws.Activate()
$ws.Range("F5").Select()
